# The author renamed the "UP" limit-type value to "FX" for every data row
# of the AF (Availability Factor) scenario sheet, i.e. column D, rows 6-229.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("AF")

$range = $ws.Range("D6:D229")
$range.Value = "FX"

# Reflect the selection the author left behind when saving the file.
$ws.Activate()
$ws.Range("D6:D229").Select()
